$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update country name labels (column A) where rank order shifted ---
$ws.Range("A31").Value = "Paises Bajos"
$ws.Range("A32").Value = "Rumania"
$ws.Range("A33").Value = "Bolivia"
$ws.Range("A48").Value = "Japon"
$ws.Range("A49").Value = "China"
$ws.Range("A86").Value = "Grecia"
$ws.Range("A87").Value = "Costa de Marfil"
$ws.Range("A109").Value = "Uganda"
$ws.Range("A110").Value = "Luxemburgo"
$ws.Range("A111").Value = "Haiti"
$ws.Range("A150").Value = "Principado de Andorra"
$ws.Range("A151").Value = "Guinea-Bisau"
$ws.Range("A152").Value = "Benin"
$ws.Range("A153").Value = "Sierra Leona"
$ws.Range("A154").Value = "Belice"
$ws.Range("A155").Value = "Burkina Faso"
$ws.Range("A156").Value = "Uruguay"
$ws.Range("A157").Value = "Letonia"
$ws.Range("A215").Value = "Islas Malvinas"
$ws.Range("A216").Value = "Montserrat"

# --- Update updated COVID numeric data (columns B-H) ---
$ws.Range("B4").Value = 7644501
$ws.Range("C4").Value = 7589
$ws.Range("D4").Value = 4861131
$ws.Range("E4").Value = 2568676
$ws.Range("G4").Value = 83
$ws.Range("H4").Value = 214694
$ws.Range("B5").Value = 6650456
$ws.Range("C5").Value = 28276
$ws.Range("D5").Value = 5621193
$ws.Range("E5").Value = 926258
$ws.Range("G5").Value = 291
$ws.Range("H5").Value = 103005
$ws.Range("B15").Value = 515571
$ws.Range("C15").Value = 12594
$ws.Range("G15").Value = 19
$ws.Range("H15").Value = 42369
$ws.Range("B21").Value = 327586
$ws.Range("C21").Value = 2257
$ws.Range("D21").Value = 232681
$ws.Range("E21").Value = 58903
$ws.Range("G21").Value = 16
$ws.Range("H21").Value = 36002
$ws.Range("B26").Value = 302542
$ws.Range("C26").Value = 971
$ws.Range("E26").Value = 31036
$ws.Range("B29").Value = 167963
$ws.Range("C29").Value = 1807
$ws.Range("D29").Value = 141660
$ws.Range("E29").Value = 16811
$ws.Range("G29").Value = 11
$ws.Range("H29").Value = 9492
$ws.Range("B31").Value = 140471
$ws.Range("C31").Value = 4579
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("G31").Value = 7
$ws.Range("H31").Value = 6461
$ws.Range("B32").Value = 137491
$ws.Range("C32").Value = 1591
$ws.Range("D32").Value = 108526
$ws.Range("E32").Value = 23917
$ws.Range("G32").Value = 45
$ws.Range("H32").Value = 5048
$ws.Range("B33").Value = 136868
$ws.Range("C33").Value = 299
$ws.Range("D33").Value = 97547
$ws.Range("E33").Value = 31220
$ws.Range("G33").Value = 28
$ws.Range("H33").Value = 8101
$ws.Range("B38").Value = 115054
$ws.Range("C38").Value = 574
$ws.Range("D38").Value = 90942
$ws.Range("E38").Value = 21968
$ws.Range("G38").Value = 10
$ws.Range("H38").Value = 2144
$ws.Range("B46").Value = 94182
$ws.Range("C46").Value = 219
$ws.Range("D46").Value = 82828
$ws.Range("E46").Value = 8052
$ws.Range("G46").Value = 9
$ws.Range("H46").Value = 3302
$ws.Range("B48").Value = 85739
$ws.Range("C48").Value = 400
$ws.Range("D48").Value = 78609
$ws.Range("E48").Value = 5531
$ws.Range("G48").Value = 2
$ws.Range("H48").Value = 1599
$ws.Range("B49").Value = 85470
$ws.Range("C49").Value = 20
$ws.Range("D49").Value = 80628
$ws.Range("E49").Value = 208
$ws.Range("H49").Value = 4634
$ws.Range("D60").Value = 57597
$ws.Range("E60").Value = 195
$ws.Range("B61").Value = 56901
$ws.Range("C61").Value = 322
$ws.Range("D61").Value = 41467
$ws.Range("E61").Value = 14059
$ws.Range("G61").Value = 9
$ws.Range("H61").Value = 1375
$ws.Range("D72").Value = 27035
$ws.Range("E72").Value = 11679
$ws.Range("G72").Value = 4
$ws.Range("H72").Value = 735
$ws.Range("B86").Value = 20142
$ws.Range("C86").Value = 300
$ws.Range("D86").Value = 9989
$ws.Range("E86").Value = 9736
$ws.Range("G86").Value = 8
$ws.Range("H86").Value = 417
$ws.Range("B87").Value = 19882
$ws.Range("D87").Value = 19449
$ws.Range("E87").Value = 313
$ws.Range("H87").Value = 120
$ws.Range("B100").Value = 12359
$ws.Range("C100").Value = 232
$ws.Range("D100").Value = 8308
$ws.Range("E100").Value = 3869
$ws.Range("G100").Value = 3
$ws.Range("H100").Value = 182
$ws.Range("B109").Value = 8965
$ws.Range("C109").Value = 157
$ws.Range("D109").Value = 5078
$ws.Range("E109").Value = 3805
$ws.Range("G109").Value = 1
$ws.Range("H109").Value = 82
$ws.Range("B110").Value = 8925
$ws.Range("C110").Value = 35
$ws.Range("D110").Value = 7793
$ws.Range("E110").Value = 1005
$ws.Range("G110").Value = 1
$ws.Range("H110").Value = 127
$ws.Range("B111").Value = 8819
$ws.Range("D111").Value = 6992
$ws.Range("E111").Value = 1598
$ws.Range("H111").Value = 229
$ws.Range("B116").Value = 7012
$ws.Range("C116").Value = 117
$ws.Range("D116").Value = 2635
$ws.Range("E116").Value = 4257
$ws.Range("B119").Value = 5845
$ws.Range("C119").Value = 36
$ws.Range("D119").Value = 5232
$ws.Range("E119").Value = 490
$ws.Range("G119").Value = 1
$ws.Range("H119").Value = 123
$ws.Range("B133").Value = 4766
$ws.Range("C133").Value = 3
$ws.Range("D133").Value = 2884
$ws.Range("E133").Value = 1801
$ws.Range("B136").Value = 4328
$ws.Range("C136").Value = 150
$ws.Range("E136").Value = 952
$ws.Range("B138").Value = 3892
$ws.Range("C138").Value = 113
$ws.Range("E138").Value = 886
$ws.Range("B143").Value = 3493
$ws.Range("C143").Value = 91
$ws.Range("E143").Value = 221
$ws.Range("B150").Value = 2370
$ws.Range("C150").Value = 260
$ws.Range("D150").Value = 1615
$ws.Range("E150").Value = 702
$ws.Range("H150").Value = 53
$ws.Range("B151").Value = 2362
$ws.Range("D151").Value = 1549
$ws.Range("E151").Value = 774
$ws.Range("H151").Value = 39
$ws.Range("B152").Value = 2357
$ws.Range("D152").Value = 1973
$ws.Range("E152").Value = 343
$ws.Range("H152").Value = 41
$ws.Range("B153").Value = 2269
$ws.Range("C153").Value = 0
$ws.Range("D153").Value = 1706
$ws.Range("E153").Value = 491
$ws.Range("G153").Value = 0
$ws.Range("H153").Value = 72
$ws.Range("B154").Value = 2196
$ws.Range("C154").Value = 65
$ws.Range("D154").Value = 1378
$ws.Range("E154").Value = 788
$ws.Range("G154").Value = 1
$ws.Range("H154").Value = 30
$ws.Range("B155").Value = 2167
$ws.Range("D155").Value = 1419
$ws.Range("E155").Value = 689
$ws.Range("H155").Value = 59
$ws.Range("B156").Value = 2145
$ws.Range("C156").Value = 0
$ws.Range("D156").Value = 1844
$ws.Range("E156").Value = 253
$ws.Range("G156").Value = 0
$ws.Range("H156").Value = 48
$ws.Range("B157").Value = 2126
$ws.Range("C157").Value = 40
$ws.Range("D157").Value = 1307
$ws.Range("E157").Value = 780
$ws.Range("G157").Value = 1
$ws.Range("H157").Value = 39
$ws.Range("B162").Value = 1847
$ws.Range("C162").Value = 23
$ws.Range("E162").Value = 456
$ws.Range("B185").Value = 344
$ws.Range("C185").Value = 2
$ws.Range("D185").Value = 316
$ws.Range("E185").Value = 4
$ws.Range("D215").Value = 13
$ws.Range("H215").Value = 0
$ws.Range("D216").Value = 12
$ws.Range("H216").Value = 1
